# Figure 2 revision #1
#  - update the auto-date placeholder text on the slide master + all 11 layouts
#  - duplicate the small eye-tracking picture, reposition both copies
#  - drop the "Speed*Occlusion Duration" connector/label pair
#  - reword/resize the "invisible" label, add two new labels
#
# All distances in the OOXML are EMUs; the COM surface reports Left/Top/
# Width/Height in points, so EMU -> pt is a divide by 12700.

$EMU = 12700

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder: "11/2/2021" -> "11/10/2021" on the slide master
#    and on every one of its custom (slide) layouts.
# ---------------------------------------------------------------------
$dateTargets = New-Object System.Collections.ArrayList
[void]$dateTargets.Add($p.SlideMaster)
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    [void]$dateTargets.Add($p.SlideMaster.CustomLayouts.Item($i))
}

foreach ($holder in $dateTargets) {
    foreach ($shp in $holder.Shapes) {
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "11/2/2021") {
                $shp.TextFrame.TextRange.Text = "11/10/2021"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1 ("Figure 2") shape edits
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

# -- 2a) Duplicate "Picture 2", move the copy in front of the original,
#        and reposition/resize both pictures.
$orig = $null
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Picture 2") {
        $orig = $shp
        break
    }
}

$dupRange = $orig.Duplicate()
$newPic = $dupRange.Item(1)
$newPic.Name = "Picture 28"
$newPic.Left = 4496837 / $EMU
$newPic.Top = 1180556 / $EMU
$newPic.Width = 549990 / $EMU
$newPic.Height = 3358604 / $EMU
$newPic.Rotation = 270
$newPic.ZOrder(1)   # msoBringToFront -> lands right before the original pic

$orig.Left = 8048290 / $EMU
$orig.Top = 1391022 / $EMU
$orig.Width = 548640 / $EMU
$orig.Height = 2875824 / $EMU

# -- 2b) Remove the dotted connector + its "Speed*Occlusion Duration" label.
$toDelete = New-Object System.Collections.ArrayList
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Straight Connector 19" -or $shp.Name -eq "TextBox 21") {
        [void]$toDelete.Add($shp)
    }
}
foreach ($shp in $toDelete) {
    $shp.Delete()
}

# -- 2c) "invisible" -> "Invisible (0.5, 0.6, 0.7s)", reposition/resize.
$invisibleBox = $null
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "TextBox 18") {
        $invisibleBox = $shp
        break
    }
}
$invisibleBox.Left = 7185142 / $EMU
$invisibleBox.Top = 3021806 / $EMU
$invisibleBox.Width = 2325637 / $EMU
$invisibleBox.TextFrame.TextRange.Text = "Invisible (0.5, 0.6, 0.7s)"

# -- 2d) Add "visible (0.5s)" label, cloned from the (now edited) box above
#        so it picks up the same body/run formatting.
$visDup = $invisibleBox.Duplicate()
$visibleBox = $visDup.Item(1)
$visibleBox.Name = "TextBox 29"
$visibleBox.Left = 3994441 / $EMU
$visibleBox.Top = 3033968 / $EMU
$visibleBox.Width = 1350050 / $EMU
$visibleBox.Height = 369332 / $EMU
$visibleBox.TextFrame.TextRange.Text = "visible (0.5s)"

# -- 2e) Add "Point of Disappearance" label (word-wraps onto two lines).
$podDup = $invisibleBox.Duplicate()
$podBox = $podDup.Item(1)
$podBox.Name = "TextBox 32"
$podBox.Left = 6737843 / $EMU
$podBox.Top = 1231296 / $EMU
$podBox.Width = 1679495 / $EMU
$podBox.Height = 646331 / $EMU
$podBox.TextFrame.WordWrap = -1
$podBox.TextFrame.TextRange.Text = "Point of Disappearance"
